$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Widen column H (Email) from 26.85546875 to 34 "characters" wide.
#    The engine adds ~0.8333333333333 padding on top of whatever we
#    set, so back that out to land exactly on 34 in the saved XML.
# ------------------------------------------------------------------
$ws.Columns.Item(8).ColumnWidth = 33.1666666666667

# ------------------------------------------------------------------
# 2) Update existing row 2 (student "Saniya Vohra") to the new
#    student "Dhruvil Patel" record.
# ------------------------------------------------------------------
$ws.Range("A2").Value = "Mr."
$ws.Range("B2").Value = "Dhruvil"
$ws.Range("C2").Value = "R"
$ws.Range("D2").Value = "Patel"
$ws.Range("E2").Value = "Male"
$ws.Range("F2").Value = "A+"
$ws.Range("G2").Value = 38456
$ws.Range("H2").Value = "panchaldhruval96@gmail.com"
$ws.Range("I2").Value = "Dhruvil@123"
$ws.Range("J2").Value = 8488887954
$ws.Range("L2").Value = 45322
$ws.Range("M2").Value = 1
$ws.Range("N2").Value = 1
$ws.Range("O2").Value = "Father"
$ws.Range("P2").Value = "Hiteshbhai"
$ws.Range("Q2").Value = "Maheshbhai"
$ws.Range("R2").Value = "Male"
$ws.Range("S2").Value = "hitesh@gmail.com"
$ws.Range("T2").Value = 8487521365
$ws.Range("U2").Value = 39844
$ws.Range("V2").Value = "MBA"
$ws.Range("W2").Value = "Job"
$ws.Range("X2").Value = "Isanpur"
$ws.Range("Y2").Value = 382443
$ws.Range("Z2").Value = "Ahmedabad"
$ws.Range("AA2").Value = "Gujarat"
$ws.Range("AB2").Value = "Isanpur"
$ws.Range("AC2").Value = 382443
$ws.Range("AD2").Value = "Ahmedabad"
$ws.Range("AE2").Value = "Gujarat"

# (G2/L2/U2 already carry the date-number-format style from the
#  original workbook, and setting .Value above does not disturb it.)

# ------------------------------------------------------------------
# 3) Append brand-new row 3 for student "Bansri More".
# ------------------------------------------------------------------
$ws.Range("A3").Value = "Mrs."
$ws.Range("B3").Value = "Bansri"
$ws.Range("C3").Value = "K"
$ws.Range("D3").Value = "More"
$ws.Range("E3").Value = "Female"
$ws.Range("F3").Value = "A+"
$ws.Range("G3").Value = 36640
$ws.Range("H3").Value = "dhruv.welinfoweb@gmail.com"
$ws.Range("I3").Value = "bansri@gmail.com"
$ws.Range("J3").Value = 8488887954
$ws.Range("L3").Value = 45322
$ws.Range("M3").Value = 1
$ws.Range("N3").Value = 1
$ws.Range("O3").Value = "Father"
$ws.Range("P3").Value = "Kamleshbhai"
$ws.Range("Q3").Value = "Baburao"
$ws.Range("R3").Value = "Male"
$ws.Range("S3").Value = "kamlesh@gmail.com"
$ws.Range("T3").Value = 8487521365
$ws.Range("U3").Value = 39844
$ws.Range("V3").Value = "MBA"
$ws.Range("W3").Value = "Job"
$ws.Range("X3").Value = "Isanpur"
$ws.Range("Y3").Value = 382443
$ws.Range("Z3").Value = "Ahmedabad"
$ws.Range("AA3").Value = "Gujarat"
$ws.Range("AB3").Value = "Isanpur"
$ws.Range("AC3").Value = 382443
$ws.Range("AD3").Value = "Ahmedabad"
$ws.Range("AE3").Value = "Gujarat"

# Row 3 is brand new, so G3/L3/U3 start out unformatted -- copy the
# date number-format from row 2's matching cells (reuses the existing
# style record instead of registering a new custom numFmt).
$ws.Range("G2").Copy()
$ws.Range("G3").PasteSpecial(-4122)
$ws.Range("L2").Copy()
$ws.Range("L3").PasteSpecial(-4122)
$ws.Range("U2").Copy()
$ws.Range("U3").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 4) Rebuild the hyperlinks: wipe the old two (their target e-mails
#    changed anyway) and re-add the six needed, in the same order the
#    authored workbook lists them (H2, S2, H3, S3, I2, I3) so the
#    relationship ids line up (rId1..rId6).
# ------------------------------------------------------------------
$ws.Range("A1:AE3").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("H2"), "mailto:panchaldhruval96@gmail.com")
$ws.Hyperlinks.Add($ws.Range("S2"), "mailto:hitesh@gmail.com")
$ws.Hyperlinks.Add($ws.Range("H3"), "mailto:dhruv.welinfoweb@gmail.com")
$ws.Hyperlinks.Add($ws.Range("S3"), "mailto:kamlesh@gmail.com")
$ws.Hyperlinks.Add($ws.Range("I2"), "mailto:Dhruvil@123")
$ws.Hyperlinks.Add($ws.Range("I3"), "mailto:bansri@gmail.com")

$ws.Range("H2").Style = "Hyperlink"
$ws.Range("S2").Style = "Hyperlink"
$ws.Range("H3").Style = "Hyperlink"
$ws.Range("S3").Style = "Hyperlink"
$ws.Range("I2").Style = "Hyperlink"
$ws.Range("I3").Style = "Hyperlink"

# ------------------------------------------------------------------
# 5) Selection lands on D7, matching the authored file.
# ------------------------------------------------------------------
$ws.Range("D7").Select()
